$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 9 ("倒進/pv/+"), shifting the
# rest of the word-list table (old rows 9-13) down to rows 10-14.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row with the new vocabulary entry: 被 (passive / 虛詞).
$ws.Range("A9").Value = "被"
$ws.Range("B9").Value = "passive"
$ws.Range("C9").Value = "虛詞"

# Highlight the new row: red font color ...
$newRow = $ws.Range("A9:C9")
$newRow.Font.Color = 255

# ... and a red medium box border drawn around the outside of A9:C9.
$newRow.Borders.Item(8).Weight = -4138
$newRow.Borders.Item(8).Color = 255
$newRow.Borders.Item(9).Weight = -4138
$newRow.Borders.Item(9).Color = 255

$ws.Range("A9").Borders.Item(7).Weight = -4138
$ws.Range("A9").Borders.Item(7).Color = 255

$ws.Range("C9").Borders.Item(10).Weight = -4138
$ws.Range("C9").Borders.Item(10).Color = 255

# Make the highlighted rows a bit taller to fit the heavier border.
$ws.Rows.Item(8).RowHeight = 16.5
$ws.Rows.Item(9).RowHeight = 16.5

# Match the author's final selection on the new row.
$newRow.Select()
